$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header/column content for the future columns
$ws.Range("D1").Value = "User_words"
$ws.Range("E1").Value = "Banned words"
$ws.Range("E2").Value = "Hitler"
$ws.Range("E3").Value = "Nazi"

# Update the active selection to match the diff (E4)
$ws.Range("E4").Select()
